$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.554.73'
$ws.Range("E2").Value = '  -1.82%  '
$ws.Range("D3").Value = '3.373.36'
$ws.Range("E3").Value = '  -2.78%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '403.12'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.81%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '132.10'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.49%  '
$ws.Range("E7").Value = '  -1.17%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.666'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.56%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.119'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -7.31%  '
$ws.Range("E11").Value = '  +1.29%  '
$ws.Range("E12").Value = '  -1.72%  '
$ws.Range("D13").Value = '3.893.48'
$ws.Range("E13").Value = '  -3.01%  '
$ws.Range("E14").Value = '  -3.08%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '19.71'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.63%  '
$ws.Range("D16").Value = '3.371.48'
$ws.Range("E16").Value = '  -2.62%  '
$ws.Range("D17").Value = '61.495.89'
$ws.Range("E17").Value = '  -1.85%  '
$ws.Range("E18").Value = '  -2.38%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.88'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.32%  '
$ws.Range("E20").Value = '  -8.13%  '
$ws.Range("E21").Value = '  -5.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '84.52'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.87%  '
$ws.Range("E23").Value = '  -0.82%  '
$ws.Range("E24").Value = '  -2.82%  '
$ws.Range("E25").Value = '  -2.89%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.78'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +10.92%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '29.34'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.13%  '
$ws.Range("E28").Value = '  +3.41%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.67'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.70'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.96%  '
$ws.Range("E31").Value = '  -2.39%  '
$ws.Range("E32").Value = '  -2.15%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.16%  '
$ws.Range("B34").Value = 'Cosmos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.29'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.57%  '
$ws.Range("B35").Value = 'InjectiveProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '41.58'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.72%  '
$ws.Range("E36").Value = '  -3.58%  '
$ws.Range("E37").Value = '  -1.29%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.998'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.04%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.40'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.56%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.93'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.77%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '138.64'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.76%  '
$ws.Range("E42").Value = '  -1.55%  '
$ws.Range("E43").Value = '  -1.30%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.292'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.88%  '
$ws.Range("E45").Value = '  +0.89%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.61'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.09%  '
$ws.Range("E47").Value = '  -1.50%  '
$ws.Range("E48").Value = '  -3.61%  '
$ws.Range("D49").Value = '2.111.98'
$ws.Range("E49").Value = '  -4.29%  '
$ws.Range("E50").Value = '  -6.37%  '
$ws.Range("E51").Value = '  +0.47%  '
